$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" value from C11 (becomes an empty inline string cell)
$ws.Range("C11").Value = ""

# Add new row 12 with the latest scraped result
# Force column A to text so the date-like string isn't auto-converted to a date serial,
# then restore the default "Normal" style so no stray number format sticks to the cell.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-03-14"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = "développement durable"
$ws.Range("C12").Value = 93
$ws.Range("D12").Value = 1
